$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (old "اسم القسم" column), shifting
# the old D/E/F ("اسم القسم"، "الدور الوظيفي"، "رقم الجوال") columns right
# by one. This is the new "الدائرة" (circle) column requested by the
# commit message ("added circle column").
$ws.Columns("D:D").Insert()

# --- Sample data row (row 2) --------------------------------------------
# Update the sample/test row first (B2's "فحص الرفع" placeholder loses its
# trailing space) so the shared-string table picks it up right after the
# untouched "رقم الجوال" entry, matching the authored template's ordering.
$ws.Range("A2").Value = 15555
$ws.Range("B2").Value = "فحص الرفع"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 59714451

# --- Header row (row 1) -------------------------------------------------
# A1 / B1 keep their original text (الرقم الوظيفي / اسم الموظف).
# Re-label the remaining headers, in left-to-right order, to match the
# new template wording.
$ws.Range("C1").Value = " المستشفى"
$ws.Range("D1").Value = " الدائرة"
$ws.Range("E1").Value = " القسم"
$ws.Range("F1").Value = " الدور الوظيفي"
# G1 keeps its original text (رقم الجوال) - only its column letter moved.

# --- Column widths ----------------------------------------------------
# (ColumnWidth assignments land 5/6 of a character wider once re-saved, so
# the values below are nudged back by 5/6 to land exactly on the template's
# target widths of 16,17,15,13,18,16,12.)
$ws.Columns("A:A").ColumnWidth = 15.1666666666667
$ws.Columns("B:B").ColumnWidth = 16.1666666666667
$ws.Columns("C:C").ColumnWidth = 14.1666666666667
$ws.Columns("D:D").ColumnWidth = 12.1666666666667
$ws.Columns("E:E").ColumnWidth = 17.1666666666667
$ws.Columns("F:F").ColumnWidth = 15.1666666666667
$ws.Columns("G:G").ColumnWidth = 11.1666666666667

# --- Selection -------------------------------------------------------------
$ws.Range("E3").Select()
